$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (serial 45406 -> 45436, i.e. 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update prices in column D for rows 29-32
$ws.Range("D29").Value = 1067
$ws.Range("D30").Value = 1265
$ws.Range("D31").Value = 2451
$ws.Range("D32").Value = 2750
